$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-27 Tuesday" "2024-08-28 Wednesday"

Replace-Text "88×20=" "56×19="
Replace-Text "47×63=" "14×21="
Replace-Text "66×60=" "72×38="
Replace-Text "19×57=" "36×82="
Replace-Text "60×85=" "11×44="

Replace-Text "80×26=" "17×85="
Replace-Text "21×11=" "14×44="
Replace-Text "17×35=" "91×40="
Replace-Text "91×60=" "62×32="
Replace-Text "61×22=" "35×21="

Replace-Text "95×85=" "90×23="
Replace-Text "34×66=" "75×35="
Replace-Text "23×88=" "85×85="
Replace-Text "91×85=" "20×96="
Replace-Text "47×12=" "98×66="

Replace-Text "75×33=" "93×98="
Replace-Text "21×41=" "12×34="
Replace-Text "64×26=" "48×61="
Replace-Text "45×71=" "20×93="
Replace-Text "99×42=" "47×23="

Replace-Text "86×85=" "43×41="
Replace-Text "37×45=" "39×97="
Replace-Text "58×19=" "13×80="
Replace-Text "92×13=" "62×21="
Replace-Text "16×46=" "28×61="
